$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (rich-text shared strings) ---
# "Volume 32   Number  12" -> "...13"  (chars 21-22, 1-based)
$ws.Range("A8").Characters(21, 2).Text = "13"

# "Report Covering the Week  3/17/2025  Through  3/23/2025"
# -> "...3/24/2025  Through  3/30/2025" (same-length date substrings)
$ws.Range("C9").Characters(27, 9).Text = "3/24/2025"
$ws.Range("C9").Characters(47, 9).Text = "3/30/2025"

# --- Helpers: change a cells literal value while forcing the exact
# cell style/number-format that the target snapshot expects, by pasting
# (format-only) from a stable donor cell that already carries that style. ---
# Style 13 = General-format "label" cells that hold literal text "0" / "***.*"
# Style 14 = "#,##0" integer cells
# Style 15 = "#,##0.0;(#,##0.0)" percent-change cells
function Set-TextCell($addr, $text, $formatSource) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
    $ws.Range($formatSource).Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = $false
}
function Set-NumberCellWithStyle($addr, $number, $formatSource) {
    $ws.Range($addr).Value = $number
    $ws.Range($formatSource).Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = $false
}

# --- Weekly crime-stat grid updates (rows 15-28) ---

# Row 15
Set-TextCell "C15" "0" "D14"
$ws.Range("M15").Value = 250
$ws.Range("N15").Value = -12.5

# Row 16
$ws.Range("C16").Value = 1
Set-NumberCellWithStyle "D16" 3 "J14"
Set-NumberCellWithStyle "E16" -66.666666666666 "K14"
$ws.Range("F16").Value = 7
$ws.Range("G16").Value = 6
$ws.Range("H16").Value = 16.666666666666
$ws.Range("J16").Value = 41
$ws.Range("K16").Value = -34.146341463414
$ws.Range("L16").Value = -22.857142857142
$ws.Range("M16").Value = -27.027027027027
$ws.Range("N16").Value = -83.832335329341

# Row 17
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 33.333333333333
$ws.Range("F17").Value = 15
$ws.Range("H17").Value = -6.25
$ws.Range("I17").Value = 53
$ws.Range("J17").Value = 40
$ws.Range("K17").Value = 32.5
$ws.Range("L17").Value = 15.217391304347
$ws.Range("M17").Value = 76.666666666666
$ws.Range("N17").Value = -7.017543859649

# Row 18
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 300
$ws.Range("F18").Value = 12
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = 20
$ws.Range("I18").Value = 31
$ws.Range("J18").Value = 29
$ws.Range("K18").Value = 6.896551724137
$ws.Range("L18").Value = -24.390243902439
$ws.Range("M18").Value = -22.5
$ws.Range("N18").Value = -88.686131386861

# Row 19
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 20
$ws.Range("E19").Value = -50
$ws.Range("F19").Value = 35
$ws.Range("G19").Value = 63
$ws.Range("H19").Value = -44.444444444444
$ws.Range("I19").Value = 112
$ws.Range("J19").Value = 195
$ws.Range("K19").Value = -42.564102564102
$ws.Range("L19").Value = -23.287671232876
$ws.Range("M19").Value = 60
$ws.Range("N19").Value = 27.272727272727

# Row 20
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = -20
$ws.Range("F20").Value = 21
$ws.Range("G20").Value = 22
$ws.Range("H20").Value = -4.545454545454
$ws.Range("I20").Value = 56
$ws.Range("J20").Value = 98
$ws.Range("K20").Value = -42.857142857142
$ws.Range("L20").Value = -31.707317073170
$ws.Range("M20").Value = 107.407407407407
$ws.Range("N20").Value = -86.138613861386

# Row 21
$ws.Range("C21").Value = 23
$ws.Range("D21").Value = 32
$ws.Range("E21").Value = -28.125
$ws.Range("F21").Value = 93
$ws.Range("G21").Value = 117
$ws.Range("H21").Value = -20.512820512820
$ws.Range("I21").Value = 286
$ws.Range("J21").Value = 405
$ws.Range("K21").Value = -29.382716049382
$ws.Range("L21").Value = -19.436619718309
$ws.Range("M21").Value = 36.842105263157
$ws.Range("N21").Value = -71.4

# Row 22
$ws.Range("L22").Value = 0

# Row 23
Set-TextCell "C23" "0" "D14"
Set-NumberCellWithStyle "D23" 3 "J14"
Set-NumberCellWithStyle "E23" -100 "K14"
$ws.Range("G23").Value = 6
$ws.Range("H23").Value = -33.333333333333
$ws.Range("J23").Value = 19
$ws.Range("K23").Value = -26.315789473684
$ws.Range("L23").Value = -26.315789473684
$ws.Range("M23").Value = 0

# Row 24
$ws.Range("C24").Value = 12
$ws.Range("D24").Value = 18
$ws.Range("E24").Value = -33.333333333333
$ws.Range("F24").Value = 67
$ws.Range("G24").Value = 82
$ws.Range("H24").Value = -18.292682926829
$ws.Range("I24").Value = 212
$ws.Range("J24").Value = 262
$ws.Range("K24").Value = -19.083969465648
$ws.Range("L24").Value = -14.17004048583
$ws.Range("M24").Value = 10.416666666666

# Row 25
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = -33.333333333333
$ws.Range("F25").Value = 12
$ws.Range("G25").Value = 36
$ws.Range("H25").Value = -66.666666666666
$ws.Range("I25").Value = 48
$ws.Range("J25").Value = 103
$ws.Range("K25").Value = -53.398058252427
$ws.Range("L25").Value = -52.475247524752

# Row 26
$ws.Range("C26").Value = 5
$ws.Range("E26").Value = -28.571428571428
$ws.Range("F26").Value = 22
$ws.Range("G26").Value = 31
$ws.Range("H26").Value = -29.032258064516
$ws.Range("I26").Value = 76
$ws.Range("J26").Value = 77
$ws.Range("K26").Value = -1.298701298701
$ws.Range("L26").Value = 8.571428571428
$ws.Range("M26").Value = -9.523809523809

# Row 27
Set-TextCell "C27" "0" "D14"

# Row 28
Set-NumberCellWithStyle "C28" 2 "J14"
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 100
$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 8
$ws.Range("H28").Value = -50
$ws.Range("I28").Value = 9
$ws.Range("J28").Value = 15
$ws.Range("K28").Value = -40
$ws.Range("L28").Value = 12.5
